$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the long note from F2 into E2 and apply the "wrap text" style
# (same style already used on G2 / B5), then clear out F2.
# NOTE: use Value2 (not Value) to avoid COM interop quirks with the
# plain Value property on this runtime.
$note = $ws.Range("F2").Value2
$ws.Range("F2").ClearContents()
$ws.Range("E2").Value2 = $note
$ws.Range("E2").WrapText = $true

# Update the active selection to E2 to match the saved view state.
$ws.Range("E2").Select()
